$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.540.04"
$ws.Range("E2").Value = "  -3.40%  "
$ws.Range("D3").Value = "'1.849.87"
$ws.Range("E3").Value = "  -3.67%  "
$ws.Range("E4").Value = "  -0.84%  "
$ws.Range("D5").Value = "'333.62"
$ws.Range("E5").Value = "  +2.19%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("D7").Value = "'0.4662"
$ws.Range("E7").Value = "  -3.21%  "
$ws.Range("D8").Value = "'0.3921"
$ws.Range("E8").Value = "  -3.49%  "
$ws.Range("D9").Value = "'46.36"
$ws.Range("E9").Value = "  -3.02%  "
$ws.Range("D10").Value = "'0.07893"
$ws.Range("E10").Value = "  -4.24%  "
$ws.Range("D11").Value = "'0.9841"
$ws.Range("E11").Value = "  -2.62%  "
$ws.Range("D12").Value = "'22.19"
$ws.Range("E12").Value = "  -5.32%  "
$ws.Range("D13").Value = "'1.924.24"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("D14").Value = "'5.853"
$ws.Range("E14").Value = "  -3.53%  "
$ws.Range("D15").Value = "'7.021"
$ws.Range("E15").Value = "  -3.14%  "
$ws.Range("D16").Value = "'0.06826"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").Value = "'1.004"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").Value = "'87.66"
$ws.Range("E18").Value = "  -4.38%  "
$ws.Range("E19").Value = "  -3.19%  "
$ws.Range("E20").Value = "  -3.24%  "
$ws.Range("D21").Value = "'1.003"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").Value = "'28.597.56"
$ws.Range("E22").Value = "  -3.27%  "
$ws.Range("D23").Value = "'5.391"
$ws.Range("E23").Value = "  -5.30%  "
$ws.Range("D24").Value = "'11.24"
$ws.Range("E24").Value = "  -5.53%  "
$ws.Range("D25").Value = "'2.200.35"
$ws.Range("E25").Value = "  +2.20%  "
$ws.Range("D26").Value = "'2.131"
$ws.Range("E26").Value = "  -2.55%  "
$ws.Range("D27").Value = "'153.73"
$ws.Range("E27").Value = "  -1.53%  "
$ws.Range("D28").Value = "'19.39"
$ws.Range("E28").Value = "  -3.20%  "
$ws.Range("D29").Value = "'6.118"
$ws.Range("E29").Value = "  -6.62%  "
$ws.Range("D30").Value = "'2.014"
$ws.Range("E30").Value = "  -4.08%  "
$ws.Range("D31").Value = "'117.24"
$ws.Range("E31").Value = "  -2.84%  "
$ws.Range("D32").Value = "'0.9746"
$ws.Range("E32").Value = "  -4.43%  "
$ws.Range("D33").Value = "'0.09423"
$ws.Range("E33").Value = "  -2.24%  "
$ws.Range("D34").Value = "'5.363"
$ws.Range("E34").Value = "  -4.68%  "
$ws.Range("D35").Value = "'3.503"
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("E36").Value = "  -2.41%  "
$ws.Range("D37").Value = "'0.06117"
$ws.Range("E37").Value = "  -3.78%  "
$ws.Range("E38").Value = "  -4.40%  "
$ws.Range("D39").Value = "'1.163"
$ws.Range("E39").Value = "  -2.24%  "
$ws.Range("B40").Value = "Frax"
$ws.Range("C40").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D40").Value = "'1.002"
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.5706"
$ws.Range("E41").Value = "  -4.23%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'7.594"
$ws.Range("E42").Value = "  -4.13%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'10.08"
$ws.Range("E43").Value = "  -6.36%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "'0.1789"
$ws.Range("E44").Value = "  -3.32%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'2.384"
$ws.Range("E45").Value = "  -3.48%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.252"
$ws.Range("E46").Value = "  -2.47%  "
$ws.Range("E47").Value = "  -3.45%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'11.79"
$ws.Range("E48").Value = "  -5.12%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.07151"
$ws.Range("E49").Value = "  -4.61%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.905"
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "'113.33"
$ws.Range("E51").Value = "  -4.68%  "
